# saved with inputs to match Figures B.4 and B.6
#
# This script reproduces (to the extent the COM surface allows) the edits
# captured in the target diff:
#   - RVC sheet: a couple of cells pick up formatting that already exists
#     elsewhere in the style table (bold label, 4-decimal number format),
#     plus two new narrow spacer columns, plus the saved selection/active
#     sheet state.
#   - FaultsPOC sheet: the "compare to figures" block (previously in
#     columns M..R) is moved left into columns K..O, and the true inputs
#     (impedance angle components B3/C3/B4/C4) are updated to the rounded
#     values used for Figures B.4/B.6, which ripples through every
#     downstream formula on the sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("RVC")
$ws2 = $wb.Worksheets.Item("FaultsPOC")

# ---------------------------------------------------------------------
# RVC sheet formatting touch-ups
# ---------------------------------------------------------------------

# O3 ("d") picks up the same bold style already used by N3/A3/etc.
$ws1.Range("O3").Font.Bold = $true

# E6/F6 pick up the 0.0000 number format already used by the rest of
# column E/F in this table.
$ws1.Range("E6").NumberFormat = "0.0000"
$ws1.Range("F6").NumberFormat = "0.0000"

# New narrow spacer columns L:M, matching the look of the FaultsPOC sheet.
$ws1.Columns("L").ColumnWidth = 1.7096354166666665
$ws1.Columns("M").ColumnWidth = 1.5299479166666665

# ---------------------------------------------------------------------
# FaultsPOC sheet: relocate the "compare to figures" block M:R -> K:O
# ---------------------------------------------------------------------

# Row 1 headers
$ws2.Range("K1").Value = $ws2.Range("M1").Value
$ws2.Range("K1").Font.Bold = $true

$ws2.Range("N1").Value = $ws2.Range("Q1").Value
$ws2.Range("N1").Font.Bold = $true

# Row 2 reference inputs (copied as literal numbers, same as before)
$ws2.Range("K2").Value = $ws2.Range("M2").Value
$ws2.Range("K2").NumberFormat = "0.00E+00"

$ws2.Range("N2").Value = $ws2.Range("Q2").Value
$ws2.Range("N2").NumberFormat = "0.00E+00"

# Row 3/4: rounded snapshots of RVC!E7/F7/J7/K7 and the new B/C inputs
$ws2.Range("K3").Value = 3.0834
$ws2.Range("L3").Value = 2.2573
$ws2.Range("N3").Value = 0.1929
$ws2.Range("O3").Value = 0.8614

$ws2.Range("K4").Value = 4.1604
$ws2.Range("L4").Value = 5.4843
$ws2.Range("N4").Value = 0.3171
$ws2.Range("O4").Value = 1.2826

# Clear out the old M/Q/R locations now that the values live in K/L/N/O.
$ws2.Range("M1").Clear()
$ws2.Range("Q1").Clear()
$ws2.Range("M2").Clear()
$ws2.Range("Q2").Clear()
$ws2.Range("M3").Clear()
$ws2.Range("Q3").Clear()
$ws2.Range("R3").Clear()
$ws2.Range("M4").Clear()
$ws2.Range("Q4").Clear()
$ws2.Range("R4").Clear()

# New narrow/wide helper columns J:M for the relocated block.
$ws2.Columns("J").ColumnWidth = 1.5299479166666665
$ws2.Columns("K").ColumnWidth = 9.619791666666666
$ws2.Columns("L").ColumnWidth = 6.346354166666667
$ws2.Columns("M").ColumnWidth = 1.7096354166666665

# ---------------------------------------------------------------------
# FaultsPOC sheet: update the true inputs to match Figures B.4/B.6
# ---------------------------------------------------------------------
$ws2.Range("B3").Value = 0.1929
$ws2.Range("C3").Value = 0.8614
$ws2.Range("B4").Value = 0.3171
$ws2.Range("C4").Value = 1.2826

# ---------------------------------------------------------------------
# Selection / active-sheet state as saved
# ---------------------------------------------------------------------

# FaultsPOC keeps a saved selection but is no longer the active tab.
$ws2.Range("N18").Select()

# RVC becomes the active sheet/tab, with its own saved selection.
$ws1.Activate()
$ws1.Range("H18").Select()
